$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update row 4 values (D4, E4, F4, G4)
$ws.Range("D4").Value = 3132666.6752300998
$ws.Range("E4").Value = 195434.692801953
$ws.Range("F4").Value = 94.403269767761202
$ws.Range("G4").Value = 5.0483090877532897

# Update row 9 values (D9, E9, F9, G9)
$ws.Range("D9").Value = 3206690.0685868501
$ws.Range("E9").Value = 295806.34573619103
$ws.Range("F9").Value = 60.6284019947052
$ws.Range("G9").Value = 5.4831590652465803

# Update the active cell selection from A18 to E17
$ws.Range("E17").Select()
